$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text so values like "27.101.77" or
# "0.07040" are not auto-converted/normalized into numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.101.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.75"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "280.93"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5292"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3524"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07040"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.34"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8163"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07809"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.882.96"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.197"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.50"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9992"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008196"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.106.61"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.120.75"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.763"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.13"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.219"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.390"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +12.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.32"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.19%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.674"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "112.54"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.394"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.377"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08914"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.176"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7462"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.896"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +9.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.418"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5326"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01882"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9719"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.18"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.323"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.214"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9990"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4599"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1371"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.439"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.71"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.531"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.64%  "
